# Updated cryptos list on Tue Jul 30 16:50:34 UTC 2024 with GitHub Actions
#
# Refreshes Price (D) / Volume(1h) (E) for every coin row from the latest
# coinranking.com snapshot, and applies the re-ranking where dogwifhat
# overtook InjectiveProtocol (rows 47/48 swap identities+new data), and
# ONDO dropped out of the top-50 in favor of Cosmos (row 51).
#
# Price values are forced to text ('-prefixed, like Excel's own "store as
# text" convention) and ClearFormats() strips the resulting "Text" number
# format again so no stray cell style is introduced; this keeps e.g.
# "577.09" / "157.90" / "0.0669" stored exactly as scraped instead of
# being auto-coerced into floating point numbers (which would mangle
# trailing zeros and introduce binary rounding noise).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.478.23"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.89%  "
$ws.Range("D3").Value = "'3.320.46"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +1.42%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'577.09"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.07%  "
$ws.Range("D6").Value = "'181.17"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.13%  "
$ws.Range("D7").Value = "'0.631"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +5.90%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  -0.69%  "
$ws.Range("E10").Value = "  +0.66%  "
$ws.Range("D11").Value = "'0.406"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.57%  "
$ws.Range("D12").Value = "'3.900.40"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.53%  "
$ws.Range("E13").Value = "  -3.26%  "
$ws.Range("D14").Value = "'26.88"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.75%  "
$ws.Range("D15").Value = "'66.547.75"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.23%  "
$ws.Range("E16").Value = "  -0.30%  "
$ws.Range("D17").Value = "'3.314.65"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.42%  "
$ws.Range("D18").Value = "'439.88"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.11%  "
$ws.Range("E19").Value = "  -0.61%  "
$ws.Range("E20").Value = "  +0.38%  "
$ws.Range("D21").Value = "'7.54"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.18%  "
$ws.Range("D22").Value = "'73.43"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.83%  "
$ws.Range("E23").Value = "  +0.23%  "
$ws.Range("E24").Value = "  +1.64%  "
$ws.Range("D25").Value = "'3.467.14"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.97%  "
$ws.Range("E26").Value = "  -1.33%  "
$ws.Range("D27").Value = "'0.198"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +5.82%  "
$ws.Range("D28").Value = "'9.12"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.77%  "
$ws.Range("E29").Value = "  -0.65%  "
$ws.Range("D30").Value = "'1.96"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.29%  "
$ws.Range("D31").Value = "'22.75"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.22%  "
$ws.Range("E32").Value = "  +0.07%  "
$ws.Range("D33").Value = "'6.77"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.15%  "
$ws.Range("E34").Value = "  -1.99%  "
$ws.Range("E35").Value = "  -1.55%  "
$ws.Range("E36").Value = "  -3.03%  "
$ws.Range("D37").Value = "'157.90"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -2.71%  "
$ws.Range("D38").Value = "'27.40"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +1.11%  "
$ws.Range("D39").Value = "'1.81"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.45%  "
$ws.Range("D40").Value = "'2.820.29"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +4.10%  "
$ws.Range("D41").Value = "'0.789"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.28%  "
$ws.Range("D42").Value = "'4.44"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.60%  "
$ws.Range("D43").Value = "'40.65"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.91%  "
$ws.Range("E44").Value = "  -1.97%  "
$ws.Range("D45").Value = "'0.0669"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("D46").Value = "'327.59"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.29%  "
$ws.Range("B47").Value = "dogwifhat"
$ws.Range("C47").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D47").Value = "'2.34"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.02%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "'23.89"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -3.02%  "
$ws.Range("E49").Value = "  +0.30%  "
$ws.Range("E50").Value = "  +3.64%  "
$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D51").Value = "'6.15"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.67%  "
